$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Histogram")

$data = @(
    @(73.69, 0.02),
    @(82.25, 0.07),
    @(90.79, 0.21),
    @(99.35, 0.22),
    @(107.89, 0.28),
    @(116.45, 0.13),
    @(124.99, 0.05),
    @(137.82, 0.02)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = $ws.Range("B1:B8")

$wb.Save()
